# The underlying commit swaps the content of ppt/theme/theme1.xml (the
# slide master's theme, originally the "Integral" palette) with
# ppt/theme/theme2.xml (the notes master's theme, originally the
# "Office Theme" palette) - i.e. after the edit, theme1.xml carries the
# "Office" colour scheme and theme2.xml carries the "Integral" colour
# scheme (font scheme / format scheme are identical between the two
# parts, so only the 12 theme colours actually change).
#
# The PowerPoint object model only exposes one writable theme colour
# scheme for this deck (SlideMaster.Theme / NotesMaster.Theme both
# resolve to the same theme part backing ppt/theme/theme1.xml), so we
# drive the swap through that: push the "Office Theme" RGB values onto
# the presentation's ThemeColorScheme, the same 12 slots (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) captured by the diff.

function Get-RGBValue($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Index -> (R, G, B) for the target ("Office Theme") palette, in the
# fixed ThemeColorScheme order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$targetRGB = @{
    1  = @(0x00, 0x00, 0x00)   # dk1
    2  = @(0xFF, 0xFF, 0xFF)   # lt1
    3  = @(0x44, 0x54, 0x6A)   # dk2
    4  = @(0xE7, 0xE6, 0xE6)   # lt2
    5  = @(0x5B, 0x9B, 0xD5)   # accent1
    6  = @(0xED, 0x7D, 0x31)   # accent2
    7  = @(0xA5, 0xA5, 0xA5)   # accent3
    8  = @(0xFF, 0xC0, 0x00)   # accent4
    9  = @(0x44, 0x72, 0xC4)   # accent5
    10 = @(0x70, 0xAD, 0x47)   # accent6
    11 = @(0x05, 0x63, 0xC1)   # hlink
    12 = @(0x95, 0x4F, 0x72)   # folHlink
}

for ($i = 1; $i -le 12; $i++) {
    $rgb = $targetRGB[$i]
    $colors.Item($i).RGB = Get-RGBValue $rgb[0] $rgb[1] $rgb[2]
}
